$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-key the data: Weapon_damage (B), Damage_dealt (D) and Hp (E) all get
#    new values while the rest of each row (Weapons/Enemy/Location/Item/Pet)
#    stays attached to the same record.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 115
$ws.Range("D2").Value = 75
$ws.Range("E2").Value = 430

$ws.Range("B3").Value = 90
$ws.Range("D3").Value = 83
$ws.Range("E3").Value = 300

$ws.Range("B4").Value = 95
$ws.Range("D4").Value = 90
$ws.Range("E4").Value = 350

$ws.Range("B5").Value = 100
$ws.Range("D5").Value = 80
$ws.Range("E5").Value = 400

$ws.Range("B6").Value = 130
$ws.Range("D6").Value = 77
$ws.Range("E6").Value = 500

$ws.Range("B7").Value = 93
$ws.Range("D7").Value = 60
$ws.Range("E7").Value = 330

$ws.Range("B8").Value = 97
$ws.Range("D8").Value = 87
$ws.Range("E8").Value = 370

$ws.Range("B9").Value = 150
$ws.Range("D9").Value = 85
$ws.Range("E9").Value = 530

$ws.Range("B10").Value = 125
$ws.Range("D10").Value = 65
$ws.Range("E10").Value = 470

$ws.Range("B11").Value = 120
$ws.Range("D11").Value = 86
$ws.Range("E11").Value = 450

# ---------------------------------------------------------------------------
# 2. Sort A1:H11 ascending by column B ("Weapon_damage"), header row excluded
#    from the reordering. This moves each record (with its formatting) to a
#    new row based on the new Weapon_damage value.
# ---------------------------------------------------------------------------
$sortRange = $ws.Range("A1:H11")
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B1:B11"))
$sortObj.SetRange($sortRange)
$sortObj.Header = 1
$sortObj.Apply()

# ---------------------------------------------------------------------------
# 3. Turn on AutoFilter for the sorted range (adds the <autoFilter> element
#    and the hidden workbook-level _FilterDatabase name).
# ---------------------------------------------------------------------------
$sortRange.AutoFilter()

$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$H`$11")
$fdName.Visible = $false

# ---------------------------------------------------------------------------
# 4. Match the reported selection left behind in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("E11").Select()
